# Insert two new price-observation rows for "Ajo" (Terminal Hortofrutícola
# Agro Chillán) into the weekly series, right before the existing row that
# used to be row 452 (date 44810). This shifts all subsequent rows down by
# two and extends the used range from A1:R518 to A1:R520.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 452; everything below (old 452..518)
# shifts down to 454..520.
$ws.Rows("452:453").Insert()

# --- New row 452 ---
$ws.Range("A452").Value = 7
$ws.Range("B452").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C452").Value = "Ñuble"
$ws.Range("D452").Value = 45180
$ws.Range("E452").Value = 16
$ws.Range("F452").Value = 100112003
$ws.Range("G452").Value = "Ajo"
$ws.Range("H452").Value = "Chino"
$ws.Range("I452").Value = "Primera"
$ws.Range("J452").Value = 40
$ws.Range("K452").Value = 21000
$ws.Range("L452").Value = 21000
$ws.Range("M452").Value = 21000
$ws.Range("N452").Value = "$/caja 10 kilos"
$ws.Range("O452").Value = "China"
$ws.Range("P452").Value = 2100
$ws.Range("Q452").Value = 10
$ws.Range("R452").Value = "Hortaliza"

# --- New row 453 ---
$ws.Range("A453").Value = 7
$ws.Range("B453").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C453").Value = "Ñuble"
$ws.Range("D453").Value = 45180
$ws.Range("E453").Value = 16
$ws.Range("F453").Value = 100112003
$ws.Range("G453").Value = "Ajo"
$ws.Range("H453").Value = "Chino"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 30
$ws.Range("K453").Value = 23000
$ws.Range("L453").Value = 23000
$ws.Range("M453").Value = 23000
$ws.Range("N453").Value = "$/malla 10 kilos"
$ws.Range("O453").Value = "China"
$ws.Range("P453").Value = 2300
$ws.Range("Q453").Value = 10
$ws.Range("R453").Value = "Hortaliza"
